# "Data driven via excel" -- repurpose the two generic sample sheets
# (Sheet1/Sheet2 holding Name/Age/Mobile_Number sample rows) into the
# two pages actually used by the Selenium E-Commerce suite:
#   Sheet1 -> ContactPage  (contact-us confirmation copy)
#   Sheet2 -> LandingPage  (store title)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename the sheets to match their new purpose.
$ws1.Name = "ContactPage"
$ws2.Name = "LandingPage"

# Wipe the old Name/Age/Mobile_Number sample data...
$ws1.Cells.Clear()
$ws2.Cells.Clear()

# ...and replace it with the strings the tests now read.
$ws1.Range("A1").Value = "CUSTOMER SERVICE - CONTACT US"
$ws1.Range("A2").Value = "Your message has been successfully sent to our team."

$ws2.Range("A1").Value = "My Store"

# LandingPage no longer holds the leftover selection on B8, and
# ContactPage (not LandingPage) becomes the active/selected tab.
$ws2.Range("A1").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
